$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows above the existing row 272, shifting the old
# rows 272-300 down to 274-302 (and extending the used range to T302).
$ws.Rows("272:273").Insert()

# Populate the two newly-inserted rows with their price-record data.
# Columns A-K (mercado/producto metadata) match every other row in this
# weekly block; columns L-T carry the new record's quality/volume/price info.

$ws.Cells.Item(272, 1).Value2 = 5
$ws.Cells.Item(272, 2).Value2 = "Macroferia Regional de Talca"
$ws.Cells.Item(272, 3).Value2 = "Maule"
$ws.Cells.Item(272, 4).Value2 = 44769
$ws.Cells.Item(272, 5).Value2 = 7
$ws.Cells.Item(272, 6).Value2 = "Fruta"
$ws.Cells.Item(272, 7).Value2 = 100101
$ws.Cells.Item(272, 8).Value2 = "Berries"
$ws.Cells.Item(272, 9).Value2 = 100101007
$ws.Cells.Item(272, 10).Value2 = "Kiwi"
$ws.Cells.Item(272, 11).Value2 = "Hayward"
$ws.Cells.Item(272, 12).Value2 = "Especial"
$ws.Cells.Item(272, 13).Value2 = 230
$ws.Cells.Item(272, 14).Value2 = 8000
$ws.Cells.Item(272, 15).Value2 = 8000
$ws.Cells.Item(272, 16).Value2 = 8000
$ws.Cells.Item(272, 17).Value2 = "`$/bandeja 18 kilos"
$ws.Cells.Item(272, 18).Value2 = "Provincia de Curicó"
$ws.Cells.Item(272, 19).Value2 = 444
$ws.Cells.Item(272, 20).Value2 = 18

$ws.Cells.Item(273, 1).Value2 = 5
$ws.Cells.Item(273, 2).Value2 = "Macroferia Regional de Talca"
$ws.Cells.Item(273, 3).Value2 = "Maule"
$ws.Cells.Item(273, 4).Value2 = 44769
$ws.Cells.Item(273, 5).Value2 = 7
$ws.Cells.Item(273, 6).Value2 = "Fruta"
$ws.Cells.Item(273, 7).Value2 = 100101
$ws.Cells.Item(273, 8).Value2 = "Berries"
$ws.Cells.Item(273, 9).Value2 = 100101007
$ws.Cells.Item(273, 10).Value2 = "Kiwi"
$ws.Cells.Item(273, 11).Value2 = "Hayward"
$ws.Cells.Item(273, 12).Value2 = "Primera"
$ws.Cells.Item(273, 13).Value2 = 200
$ws.Cells.Item(273, 14).Value2 = 7000
$ws.Cells.Item(273, 15).Value2 = 7000
$ws.Cells.Item(273, 16).Value2 = 7000
$ws.Cells.Item(273, 17).Value2 = "`$/bandeja 18 kilos"
$ws.Cells.Item(273, 18).Value2 = "Provincia de Curicó"
$ws.Cells.Item(273, 19).Value2 = 389
$ws.Cells.Item(273, 20).Value2 = 18
